# Updated symbol list on Wed Dec 14 11:52:28 UTC 2022 with GitHub Actions
# Refresh the "Price" column (D) with newly scraped values and fix a
# couple of "Volume(1h)" (E) labels that no longer carry the
# Best/Worst-in-24h badge text.
#
# The Price column stores numeric-looking values as text, so a leading
# apostrophe is used to force the COM layer to keep them as text instead
# of silently re-interpreting them as numbers (which would also introduce
# binary floating point noise, e.g. 271.76 -> 271.75999999999999).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $ws.Range($cellRef).Value = "'" + $value
}

Set-TextValue "D2"  "271.76"
Set-TextValue "D3"  "22.85"
Set-TextValue "D4"  "6.348"
Set-TextValue "D5"  "0.06213"
Set-TextValue "D7"  "6.691"
Set-TextValue "D9"  "0.8348"
Set-TextValue "D11" "0.1638"
Set-TextValue "D12" "0.08316"
Set-TextValue "D13" "0.03355"
Set-TextValue "D14" "0.03103"
Set-TextValue "D15" "0.09323"
Set-TextValue "D16" "3.914"
Set-TextValue "D17" "0.001641"
Set-TextValue "D19" "0.006242"

Set-TextValue "D20" "0.005566"
$ws.Range("E20").Value = "19HotbitTokenHTB"

Set-TextValue "D21" "0.001088"
Set-TextValue "D23" "3.728"
Set-TextValue "D25" "0.3403"
Set-TextValue "D41" "0.007023"
Set-TextValue "D42" "0.1166"

Set-TextValue "D43" "0.003342"
$ws.Range("E43").Value = "42CEJICEJIWorstin24h"

Set-TextValue "D44" "0.01241"
Set-TextValue "D45" "0.00006260"
Set-TextValue "D47" "0.9002"
Set-TextValue "D48" "0.04447"
Set-TextValue "D49" "0.00002300"
Set-TextValue "D50" "0.01240"
